$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

Set-TextValue 'D2' '69.221.33'
Set-TextValue 'E2' '  +3.30%  '
Set-TextValue 'D3' '3.594.24'
Set-TextValue 'E3' '  +3.29%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.11%  '
Set-TextValue 'D5' '625.19'
Set-TextValue 'E5' '  +3.06%  '
Set-TextValue 'D6' '158.05'
Set-TextValue 'E6' '  +6.47%  '
Set-TextValue 'D7' '3.591.46'
Set-TextValue 'E7' '  +3.21%  '
Set-TextValue 'E8' '  -0.04%  '
Set-TextValue 'E9' '  +2.78%  '
Set-TextValue 'D10' '0.150'
Set-TextValue 'E10' '  +9.21%  '
Set-TextValue 'D11' '7.43'
Set-TextValue 'E11' '  +7.79%  '
Set-TextValue 'D12' '0.442'
Set-TextValue 'E12' '  +4.95%  '
Set-TextValue 'D13' '0.0000226'
Set-TextValue 'E13' '  +6.67%  '
Set-TextValue 'D14' '33.62'
Set-TextValue 'E14' '  +7.51%  '
Set-TextValue 'D15' '4.205.24'
Set-TextValue 'E15' '  +3.40%  '
Set-TextValue 'D16' '3.600.00'
Set-TextValue 'E16' '  +3.51%  '
Set-TextValue 'D17' '69.347.68'
Set-TextValue 'E17' '  +3.66%  '
Set-TextValue 'E18' '  +0.79%  '
Set-TextValue 'D19' '6.81'
Set-TextValue 'E19' '  +5.96%  '
Set-TextValue 'D20' '16.13'
Set-TextValue 'E20' '  +7.50%  '
Set-TextValue 'E21' '  +12.79%  '
Set-TextValue 'D22' '462.27'
Set-TextValue 'E22' '  +4.14%  '
Set-TextValue 'D23' '0.645'
Set-TextValue 'E23' '  +3.47%  '
Set-TextValue 'D24' '78.76'
Set-TextValue 'E24' '  +2.00%  '
Set-TextValue 'D25' '0.0000135'
Set-TextValue 'E25' '  +8.64%  '
Set-TextValue 'D26' '10.70'
Set-TextValue 'E26' '  +5.88%  '
Set-TextValue 'D27' '3.736.19'
Set-TextValue 'E27' '  +3.29%  '
Set-TextValue 'D28' '0.999'
Set-TextValue 'E28' '  -0.24%  '
Set-TextValue 'D29' '9.26'
Set-TextValue 'E29' '  +11.62%  '
Set-TextValue 'D30' '2.64'
Set-TextValue 'E30' '  +4.42%  '
Set-TextValue 'D31' '1.71'
Set-TextValue 'E31' '  +9.15%  '
Set-TextValue 'D32' '0.173'
Set-TextValue 'E32' '  +5.47%  '
Set-TextValue 'E33' '  -0.02%  '
Set-TextValue 'D34' '6.55'
Set-TextValue 'E34' '  +7.06%  '
Set-TextValue 'B35' 'EthereumClassic'
Set-TextValue 'C35' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D35' '26.48'
Set-TextValue 'E35' '  +3.50%  '
Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '1.94'
Set-TextValue 'E36' '  +4.63%  '
Set-TextValue 'D37' '3.588.21'
Set-TextValue 'E37' '  +3.70%  '
Set-TextValue 'D38' '8.42'
Set-TextValue 'E38' '  +5.74%  '
Set-TextValue 'D39' '2.40'
Set-TextValue 'E39' '  +9.63%  '
Set-TextValue 'E40' '  +0.00%  '
Set-TextValue 'B41' 'Monero'
Set-TextValue 'C41' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D41' '179.79'
Set-TextValue 'E41' '  +5.30%  '
Set-TextValue 'B42' 'Hedera'
Set-TextValue 'C42' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D42' '0.0925'
Set-TextValue 'E42' '  +6.53%  '
Set-TextValue 'D43' '1.00'
Set-TextValue 'E43' '  +0.13%  '
Set-TextValue 'D44' '5.72'
Set-TextValue 'E44' '  +5.39%  '
Set-TextValue 'D45' '31.47'
Set-TextValue 'E45' '  +21.79%  '
Set-TextValue 'D46' '0.911'
Set-TextValue 'E46' '  +3.44%  '
Set-TextValue 'E47' '  +11.54%  '
Set-TextValue 'D48' '45.95'
Set-TextValue 'E48' '  +0.64%  '
Set-TextValue 'D49' '2.75'
Set-TextValue 'E49' '  +9.78%  '
Set-TextValue 'D50' '7.82'
Set-TextValue 'E50' '  +3.61%  '
Set-TextValue 'E51' '  +8.75%  '
